# Apply cryptos list update (Tue Sep 26 18:44:10 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.162.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "'1.585.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("D5").Value = "'211.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").Value = "'1.808.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "'1.578.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "'4.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "'63.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").Value = "'26.188.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "'214.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("D25").Value = "'144.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("D33").Value = "'1.408.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.02%  "
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.589"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.19%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'1.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("E38").Value = "  -1.42%  "
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("E40").Value = "  +4.00%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "'0.938"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -15.91%  "
$ws.Range("D43").Value = "'0.765"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").Value = "'1.720.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").Value = "'60.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.82%  "
$ws.Range("D47").Value = "'85.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.97%  "
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("D49").Value = "'0.0499"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("D50").Value = "'0.0971"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.09%  "
